# Apply the "Added evaluation for some constellations." edit to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Header row (row 1) changes ----
$ws.Cells.Item(1, 9).Value = "multiple det."     # I1: was "status"
$ws.Cells.Item(1, 10).Value = "solved with"      # J1: was "reason"
$ws.Cells.Item(1, 11).Value = "abbreviations"    # K1: new column

# ---- Row 2 ----
$ws.Cells.Item(2, 4).Value = 8                   # D2: 10 -> 8
$ws.Cells.Item(2, 9).Value = 1                   # I2
$ws.Cells.Item(2, 10).Value = "GreLum, LowDiff"  # J2
$ws.Cells.Item(2, 11).Value = "LowDiff = take set with lowest triangle difference"  # K2

# ---- Row 3 ----
$ws.Cells.Item(3, 11).Value = "GreLum = take set with greatest luminosity"  # K3

# ---- Row 4 ----
$ws.Cells.Item(4, 9).Value = 0                   # I4

# ---- Row 5 ----
$ws.Cells.Item(5, 9).Value = 0                   # I5

# ---- Row 6 ----
$ws.Cells.Item(6, 9).Value = 1                   # I6
$ws.Cells.Item(6, 10).Value = "GreLum, LowDiff"  # J6

# ---- Row 7 ----
$ws.Cells.Item(7, 2).Value = 0.02                # B7
$ws.Cells.Item(7, 3).Value = "3 or 5"            # C7 (text; special style below)
$ws.Cells.Item(7, 4).Value = 5                   # D7
$ws.Cells.Item(7, 5).Value = 10                  # E7
$ws.Cells.Item(7, 6).Value = 10                  # F7
$ws.Cells.Item(7, 9).Value = 0                   # I7

# ---- Row 8 ----
$ws.Cells.Item(8, 2).Value = 0.1                 # B8
$ws.Cells.Item(8, 3).Value = 5                   # C8
$ws.Cells.Item(8, 4).Value = 5                   # D8
$ws.Cells.Item(8, 5).Value = 10                  # E8
$ws.Cells.Item(8, 6).Value = 9                   # F8
$ws.Cells.Item(8, 9).Value = 2                   # I8
$ws.Cells.Item(8, 10).Value = "GreLum, LowDiff"  # J8

# ---- Row 9 ----
$ws.Cells.Item(9, 2).Value = 0.05                # B9
$ws.Cells.Item(9, 3).Value = 3                   # C9
$ws.Cells.Item(9, 4).Value = 5                   # D9
$ws.Cells.Item(9, 5).Value = 10                  # E9
$ws.Cells.Item(9, 6).Value = 10                  # F9
$ws.Cells.Item(9, 9).Value = 0                   # I9

# ---- Row 10 ----
$ws.Cells.Item(10, 2).Value = 0.1                # B10
$ws.Cells.Item(10, 3).Value = 2                  # C10
$ws.Cells.Item(10, 4).Value = 10                 # D10
$ws.Cells.Item(10, 5).Value = 10                 # E10
$ws.Cells.Item(10, 6).Value = 3                  # F10
$ws.Cells.Item(10, 9).Value = 2                  # I10
$ws.Cells.Item(10, 10).Value = "GreLum, LowDiff. Remove some stars from model for better results?"  # J10

# ---- Row 11 ----
$ws.Cells.Item(11, 2).Value = 0.02               # B11
$ws.Cells.Item(11, 3).Value = 5                  # C11
$ws.Cells.Item(11, 4).Value = 3                  # D11
$ws.Cells.Item(11, 5).Value = 10                 # E11
$ws.Cells.Item(11, 6).Value = 10                 # F11
$ws.Cells.Item(11, 9).Value = 0                  # I11

# ---- Row 12 ----
$ws.Cells.Item(12, 2).Value = 0.1                # B12
$ws.Cells.Item(12, 3).Value = 3                  # C12
$ws.Cells.Item(12, 4).Value = 5                  # D12
$ws.Cells.Item(12, 5).Value = 7                  # E12
$ws.Cells.Item(12, 6).Value = 8                  # F12
$ws.Cells.Item(12, 9).Value = 0                  # I12

# ---- Row 13 ----
$ws.Cells.Item(13, 2).Value = 0.1                # B13
$ws.Cells.Item(13, 3).Value = 3                  # C13
$ws.Cells.Item(13, 4).Value = 5                  # D13
$ws.Cells.Item(13, 5).Value = 10                 # E13
$ws.Cells.Item(13, 6).Value = 30                 # F13
$ws.Cells.Item(13, 7).Value = 30                 # G13: 10 -> 30
$ws.Cells.Item(13, 9).Value = 2                  # I13
$ws.Cells.Item(13, 10).Value = "GreLum, LowDiff" # J13

# ---- Row 14: unchanged ----

# ---- Row 15 ----
$ws.Cells.Item(15, 2).Value = 0.1                # B15
$ws.Cells.Item(15, 3).Value = 2                  # C15
$ws.Cells.Item(15, 4).Value = 5                  # D15
$ws.Cells.Item(15, 5).Value = 10                 # E15
$ws.Cells.Item(15, 6).Value = 10                 # F15
$ws.Cells.Item(15, 9).Value = 0                  # I15

# ---- Row 16: unchanged ----

# ---- Row 17 ----
$ws.Cells.Item(17, 2).Value = 0.05               # B17
$ws.Cells.Item(17, 3).Value = 2                  # C17
$ws.Cells.Item(17, 4).Value = 5                  # D17
$ws.Cells.Item(17, 5).Value = 10                 # E17
$ws.Cells.Item(17, 6).Value = 5                  # F17
$ws.Cells.Item(17, 9).Value = 0                  # I17

# ---- Row 18 ----
$ws.Cells.Item(18, 2).Value = 0.1                # B18
$ws.Cells.Item(18, 3).Value = 3                  # C18
$ws.Cells.Item(18, 4).Value = 10                 # D18
$ws.Cells.Item(18, 5).Value = 10                 # E18
$ws.Cells.Item(18, 6).Value = 8                  # F18
$ws.Cells.Item(18, 9).Value = 1                  # I18
$ws.Cells.Item(18, 10).Value = "GreLum or LowDiff"  # J18

# ---- Row 19 ----
$ws.Cells.Item(19, 2).Value = 0.02               # B19
$ws.Cells.Item(19, 3).Value = 2                  # C19
$ws.Cells.Item(19, 4).Value = 7                  # D19
$ws.Cells.Item(19, 5).Value = 10                 # E19
$ws.Cells.Item(19, 6).Value = 10                 # F19
$ws.Cells.Item(19, 9).Value = 0                  # I19

# ---- Row 20 ----
$ws.Cells.Item(20, 2).Value = 0.1                # B20
$ws.Cells.Item(20, 3).Value = 3                  # C20
$ws.Cells.Item(20, 4).Value = 7                  # D20
$ws.Cells.Item(20, 5).Value = 10                 # E20
$ws.Cells.Item(20, 6).Value = 2                  # F20
$ws.Cells.Item(20, 9).Value = 0                  # I20

# ---- Column widths ----
$ws.Range("A:E").ColumnWidth = 13.7109375
$ws.Range("F:F").ColumnWidth = 21.5703125
$ws.Range("G:H").ColumnWidth = 13.7109375
$ws.Range("I:I").ColumnWidth = 17.85546875
$ws.Range("J:J").ColumnWidth = 34.42578125
$ws.Range("K:K").ColumnWidth = 20.7109375

# ---- Cell alignment: K column new cells left as general (default), K1 centered like other headers ----
$ws.Cells.Item(1, 11).HorizontalAlignment = -4108   # xlCenter

# ---- Special style for C7: centered + number format "d-mmm" (numFmtId 16) ----
$ws.Cells.Item(7, 3).HorizontalAlignment = -4108    # xlCenter
$ws.Cells.Item(7, 3).NumberFormat = "d-mmm"

# ---- View: scroll position & selection ----
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("I20").Select()
